# Update the "取得日時" (retrieved datetime) timestamp in column A for all
# data rows (2-9) on the "ランサーズ" sheet to reflect the new run time.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-11-24 12:51:04"

for ($row = 2; $row -le 9; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
